# Refresh the Entsoe "Actual Production - Wind" sheet with the latest pull:
# shift every timestamp forward by 12 days (new data day) and update the
# production values (MW) for the hours that now include Horeco's output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp

# New "Actual Production (MW)" readings for the rows that changed (1-based
# sheet row numbers). Rows not listed here keep their previous value (0).
$newProduction = @{
    2 = 69;  3 = 50;  4 = 34;  5 = 25;  6 = 9;   7 = 2;   8 = 4;
    9 = 5;   10 = 5;  11 = 0;  12 = 0;  13 = 0;  14 = 2;  15 = 0;
    16 = 0;  17 = 2;  18 = 2;  19 = 1;  20 = 3;  21 = 5;  22 = 11;
    23 = 17; 24 = 25; 25 = 29; 26 = 34; 27 = 48; 28 = 59; 29 = 0;
    30 = 0;  31 = 0;  32 = 0;  33 = 0;  34 = 0;  35 = 0;  36 = 0;
    37 = 0;  38 = 0;  39 = 0
}

for ($row = 2; $row -le $lastRow; $row++) {
    $tsCell = $ws.Cells.Item($row, 1)
    $tsCell.Value2 = $tsCell.Value2 + 12

    if ($newProduction.ContainsKey($row)) {
        $ws.Cells.Item($row, 2).Value2 = $newProduction[$row]
    }
}
